$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (column C) date value for rows 2-9 from 45221 to 45224
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value = 45224
}
